$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.928.58'
$ws.Range("E2").Value = '  -2.48%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.165.66'
$ws.Range("E3").Value = '  -7.74%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '566.86'
$ws.Range("E5").Value = '  -3.18%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.71'
$ws.Range("E6").Value = '  -4.76%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.619'
$ws.Range("E7").Value = '  -1.28%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.164.17'
$ws.Range("E9").Value = '  -7.78%  '

$ws.Range("E10").Value = '  -6.01%  '

$ws.Range("E11").Value = '  -5.56%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.396'
$ws.Range("E12").Value = '  -5.20%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.716.38'
$ws.Range("E13").Value = '  -7.79%  '

$ws.Range("E14").Value = '  +0.79%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.35'
$ws.Range("E15").Value = '  -7.93%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.795.59'
$ws.Range("E16").Value = '  -2.73%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000163'
$ws.Range("E17").Value = '  -6.07%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.163.79'
$ws.Range("E18").Value = '  -7.94%  '

$ws.Range("E19").Value = '  -3.65%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.97'
$ws.Range("E20").Value = '  -6.52%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '356.67'
$ws.Range("E21").Value = '  -4.13%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.27'
$ws.Range("E22").Value = '  -5.27%  '

$ws.Range("E23").Value = '  +0.35%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.74'
$ws.Range("E24").Value = '  -5.86%  '

$ws.Range("B25").Value = 'Polygon'
$ws.Range("C25").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.502'
$ws.Range("E25").Value = '  -6.80%  '

$ws.Range("B26").Value = 'PEPE'
$ws.Range("C26").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000118'
$ws.Range("E26").Value = '  -7.85%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.69'
$ws.Range("E27").Value = '  -2.30%  '

$ws.Range("E28").Value = '  -1.93%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.50%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.997'
$ws.Range("E30").Value = '  -0.22%  '

$ws.Range("E31").Value = '  -5.11%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.40'
$ws.Range("E32").Value = '  -7.84%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '22.04'
$ws.Range("E33").Value = '  -6.70%  '

$ws.Range("B34").Value = 'Aptos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.69'
$ws.Range("E34").Value = '  -5.65%  '

$ws.Range("B35").Value = 'Fetch.AI'
$ws.Range("C35").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.21'
$ws.Range("E35").Value = '  -5.65%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.45'
$ws.Range("E36").Value = '  -7.55%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '153.89'
$ws.Range("E37").Value = '  -5.57%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.833'
$ws.Range("E38").Value = '  -5.25%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '26.20'
$ws.Range("E39").Value = '  -6.14%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.75'
$ws.Range("E40").Value = '  -3.02%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.55'
$ws.Range("E41").Value = '  -3.78%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.664.77'
$ws.Range("E42").Value = '  -2.61%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.20'
$ws.Range("E43").Value = '  -6.46%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.03'
$ws.Range("E44").Value = '  -5.89%  '

$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '24.48'
$ws.Range("E45").Value = '  -4.04%  '

$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0656'
$ws.Range("E46").Value = '  -5.55%  '

$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '39.20'
$ws.Range("E47").Value = '  -2.41%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '324.46'
$ws.Range("E48").Value = '  -3.30%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0273'
$ws.Range("E49").Value = '  -5.25%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.103'
$ws.Range("E50").Value = '  -2.17%  '

$ws.Range("E51").Value = '  -0.16%  '
